# Updated cryptos list: refresh Price (D) and Volume(1h) (E) columns
# to match the latest scrape, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "57.395.50"
$ws.Range("E2").Value = "  -3.93%  "
# Row 3
$ws.Range("D3").Value = "2.935.85"
$ws.Range("E3").Value = "  -0.27%  "
# Row 4
$ws.Range("E4").Value = "  -0.01%  "
# Row 5
$ws.Range("D5").Value = "'551.85"
$ws.Range("E5").Value = "  -3.23%  "
# Row 6
$ws.Range("D6").Value = "'130.55"
$ws.Range("E6").Value = "  +6.50%  "
# Row 7
$ws.Range("E7").Value = "  -0.03%  "
# Row 8
$ws.Range("D8").Value = "'0.513"
$ws.Range("E8").Value = "  +3.35%  "
# Row 9
$ws.Range("D9").Value = "2.924.51"
$ws.Range("E9").Value = "  -0.42%  "
# Row 10
$ws.Range("E10").Value = "  -2.48%  "
# Row 11
$ws.Range("D11").Value = "'4.78"
$ws.Range("E11").Value = "  -4.88%  "
# Row 12
$ws.Range("E12").Value = "  +2.25%  "
# Row 13
$ws.Range("D13").Value = "'0.0000221"
$ws.Range("E13").Value = "  +0.57%  "
# Row 14
$ws.Range("D14").Value = "'32.50"
$ws.Range("E14").Value = "  +0.88%  "
# Row 15
$ws.Range("E15").Value = "  +1.44%  "
# Row 16
$ws.Range("D16").Value = "3.421.74"
$ws.Range("E16").Value = "  -0.66%  "
# Row 17
$ws.Range("D17").Value = "'6.74"
$ws.Range("E17").Value = "  +10.03%  "
# Row 18
$ws.Range("D18").Value = "2.932.17"
$ws.Range("E18").Value = "  -0.93%  "
# Row 19
$ws.Range("D19").Value = "57.451.19"
$ws.Range("E19").Value = "  -4.10%  "
# Row 20
$ws.Range("D20").Value = "'415.49"
$ws.Range("E20").Value = "  -3.37%  "
# Row 21
$ws.Range("D21").Value = "'13.04"
$ws.Range("E21").Value = "  +0.69%  "
# Row 22
$ws.Range("E22").Value = "  +3.47%  "
# Row 23
$ws.Range("D23").Value = "'6.92"
$ws.Range("E23").Value = "  +0.49%  "
# Row 24
$ws.Range("D24").Value = "'12.85"
$ws.Range("E24").Value = "  +1.25%  "
# Row 25
$ws.Range("D25").Value = "'78.79"
$ws.Range("E25").Value = "  +0.76%  "
# Row 26
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.09%  "
# Row 27
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.09%  "
# Row 28
$ws.Range("E28").Value = "  -1.11%  "
# Row 29
$ws.Range("E29").Value = "  +5.27%  "
# Row 30
$ws.Range("D30").Value = "'1.98"
$ws.Range("E30").Value = "  +6.46%  "
# Row 31
$ws.Range("D31").Value = "'6.10"
$ws.Range("E31").Value = "  +0.85%  "
# Row 32
$ws.Range("E32").Value = "  +12.27%  "
# Row 33
$ws.Range("D33").Value = "'24.97"
$ws.Range("E33").Value = "  -0.32%  "
# Row 34
$ws.Range("D34").Value = "'5.60"
$ws.Range("E34").Value = "  +1.54%  "
# Row 35
$ws.Range("E35").Value = "  -1.76%  "
# Row 36
$ws.Range("D36").Value = "'0.929"
$ws.Range("E36").Value = "  -0.89%  "
# Row 37
$ws.Range("D37").Value = "'48.50"
$ws.Range("E37").Value = "  -1.44%  "
# Row 38
$ws.Range("D38").Value = "0.0₃0676"
$ws.Range("E38").Value = "  +5.37%  "
# Row 39
$ws.Range("D39").Value = "'8.37"
$ws.Range("E39").Value = "  +6.89%  "
# Row 40
$ws.Range("E40").Value = "  +6.00%  "
# Row 41
$ws.Range("E41").Value = "  -1.93%  "
# Row 42
$ws.Range("D42").Value = "'0.107"
$ws.Range("E42").Value = "  +0.49%  "
# Row 43
$ws.Range("D43").Value = "'375.93"
$ws.Range("E43").Value = "  +0.57%  "
# Row 44
$ws.Range("D44").Value = "2.632.86"
$ws.Range("E44").Value = "  +1.05%  "
# Row 46
$ws.Range("D46").Value = "'0.238"
$ws.Range("E46").Value = "  +1.98%  "
# Row 47
$ws.Range("D47").Value = "'121.40"
$ws.Range("E47").Value = "  +2.83%  "
# Row 48
$ws.Range("E48").Value = "  +3.17%  "
# Row 49
$ws.Range("D49").Value = "'1.96"
$ws.Range("E49").Value = "  +1.43%  "
# Row 50
$ws.Range("D50").Value = "'23.18"
$ws.Range("E50").Value = "  +0.60%  "
# Row 51
$ws.Range("D51").Value = "'1.98"
$ws.Range("E51").Value = "  +1.60%  "
